# Updates FFXIV leveling-profit calculations (H:N) across multiple job sheets
# per the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 309.83334
$ws.Range("I28").Value = 309.83334
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 309.83334
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 175.16666
$ws.Range("N28").ClearContents()
# Row 62
$ws.Range("H62").Value = 3159.9092
$ws.Range("I62").Value = 2257.1428
$ws.Range("J62").Value = 4739.75
$ws.Range("K62").Value = 2257.1428
$ws.Range("L62").Value = 4739.75
$ws.Range("M62").Value = -1633.1428
$ws.Range("N62").Value = -5987.75
# Row 65
$ws.Range("H65").Value = 3159.9092
$ws.Range("I65").Value = 2257.1428
$ws.Range("J65").Value = 4739.75
$ws.Range("K65").Value = 11285.714
$ws.Range("L65").Value = 23698.75
$ws.Range("M65").Value = -8165.714
$ws.Range("N65").Value = -29938.75
# Row 103
$ws.Range("H103").Value = 472.65
$ws.Range("I103").Value = 465.69232
$ws.Range("J103").Value = 485.57144
$ws.Range("K103").Value = 1397.07696
$ws.Range("L103").Value = 1456.71432
$ws.Range("M103").Value = -811.0769599999999
$ws.Range("N103").Value = -2628.71432
# Row 132
$ws.Range("H132").Value = 2586.146
$ws.Range("I132").Value = 2334.039
$ws.Range("J132").Value = 4203.8335
$ws.Range("K132").Value = 7002.117
$ws.Range("L132").Value = 12611.5005
$ws.Range("M132").Value = -4472.117
$ws.Range("N132").Value = -17671.5005

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1050
$ws.Range("I2").Value = 809.375
$ws.Range("K2").Value = 809.375
$ws.Range("M2").Value = -696.375
# Row 45
$ws.Range("H45").Value = 1543.3158
$ws.Range("I45").Value = 1601.7858
$ws.Range("K45").Value = 1601.7858
$ws.Range("M45").Value = -1224.7858
# Row 110
$ws.Range("H110").Value = 1180.3889
$ws.Range("I110").Value = 976.5333000000001
$ws.Range("J110").Value = 2199.6667
$ws.Range("K110").Value = 976.5333000000001
$ws.Range("L110").Value = 2199.6667
$ws.Range("M110").Value = 1068.4667
$ws.Range("N110").Value = -6289.6667
# Row 116
$ws.Range("H116").Value = 1050
$ws.Range("I116").Value = 809.375
$ws.Range("K116").Value = 809.375
$ws.Range("M116").Value = 1484.625

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1050
$ws.Range("I3").Value = 809.375
$ws.Range("K3").Value = 809.375
$ws.Range("M3").Value = -695.375
# Row 107
$ws.Range("H107").Value = 35653
$ws.Range("I107").Value = 51489.5
$ws.Range("J107").Value = 3980
$ws.Range("K107").Value = 51489.5
$ws.Range("L107").Value = 3980
$ws.Range("M107").Value = -49569.5
$ws.Range("N107").Value = -7820

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 105
$ws.Range("H105").Value = 918.1667
$ws.Range("I105").Value = 918.1667
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 918.1667
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 828.8333
$ws.Range("N105").ClearContents()
# Row 107
$ws.Range("H107").Value = 1075.8422
$ws.Range("I107").Value = 1101.3636
$ws.Range("J107").Value = 1040.75
$ws.Range("K107").Value = 1101.3636
$ws.Range("L107").Value = 1040.75
$ws.Range("M107").Value = 818.6364000000001
$ws.Range("N107").Value = -4880.75
# Row 108
$ws.Range("H108").Value = 26000
$ws.Range("J108").Value = 26000
$ws.Range("L108").Value = 26000
$ws.Range("N108").Value = -33680

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 765.1836499999999
$ws.Range("I5").Value = 451.1154
$ws.Range("J5").Value = 1120.2174
$ws.Range("K5").Value = 1353.3462
$ws.Range("L5").Value = 3360.6522
$ws.Range("M5").Value = -1241.3462
$ws.Range("N5").Value = -3584.6522
# Row 105
$ws.Range("H105").Value = 24444.4
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 24444.4
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 73333.20000000001
$ws.Range("M105").ClearContents()
$ws.Range("N105").Value = -78575.20000000001
# Row 122
$ws.Range("H122").Value = 803.21875
$ws.Range("I122").Value = 442.3684
$ws.Range("J122").Value = 1330.6154
$ws.Range("K122").Value = 3981.3156
$ws.Range("L122").Value = 11975.5386
$ws.Range("M122").Value = -1531.3156
$ws.Range("N122").Value = -16875.5386
# Row 135
$ws.Range("H135").Value = 765.1836499999999
$ws.Range("I135").Value = 451.1154
$ws.Range("J135").Value = 1120.2174
$ws.Range("K135").Value = 4060.0386
$ws.Range("L135").Value = 10081.9566
$ws.Range("M135").Value = -1525.0386
$ws.Range("N135").Value = -15151.9566

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 1114.4231
$ws.Range("I97").Value = 1083.4736
$ws.Range("J97").Value = 1198.4286
$ws.Range("K97").Value = 1083.4736
$ws.Range("L97").Value = 1198.4286
$ws.Range("M97").Value = -587.4736
$ws.Range("N97").Value = -2190.4286
# Row 102
$ws.Range("H102").Value = 1573.25
$ws.Range("I102").Value = 1554.6666
$ws.Range("J102").Value = 1597.1428
$ws.Range("K102").Value = 1554.6666
$ws.Range("L102").Value = 1597.1428
$ws.Range("M102").Value = 67.33339999999998
$ws.Range("N102").Value = -4841.1428
# Row 107
$ws.Range("H107").Value = 558.2353000000001
$ws.Range("I107").Value = 599.6667
$ws.Range("J107").Value = 247.5
$ws.Range("K107").Value = 599.6667
$ws.Range("L107").Value = 247.5
$ws.Range("M107").Value = 1320.3333
$ws.Range("N107").Value = -4087.5
# Row 113
$ws.Range("H113").Value = 1365.5454
$ws.Range("I113").Value = 1231.7142
$ws.Range("J113").Value = 1599.75
$ws.Range("K113").Value = 1231.7142
$ws.Range("L113").Value = 1599.75
$ws.Range("M113").Value = 938.2858000000001
$ws.Range("N113").Value = -5939.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 595.25
$ws.Range("I16").Value = 595.25
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 595.25
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -425.25
$ws.Range("N16").ClearContents()
# Row 61
$ws.Range("H61").Value = 2304.5
$ws.Range("I61").Value = 2282.2222
$ws.Range("J61").Value = 2505
$ws.Range("K61").Value = 2282.2222
$ws.Range("L61").Value = 2505
$ws.Range("M61").Value = -2080.2222
$ws.Range("N61").Value = -2909
# Row 113
$ws.Range("H113").Value = 2304.5
$ws.Range("I113").Value = 2282.2222
$ws.Range("J113").Value = 2505
$ws.Range("K113").Value = 2282.2222
$ws.Range("L113").Value = 2505
$ws.Range("M113").Value = -112.2222000000002
$ws.Range("N113").Value = -6845
# Row 136
$ws.Range("H136").Value = 13504.583
$ws.Range("I136").Value = 17100.75
$ws.Range("K136").Value = 51302.25
$ws.Range("M136").Value = -48752.25

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 2129
$ws.Range("I107").Value = 3810
$ws.Range("K107").Value = 11430
$ws.Range("M107").Value = -9510
# Row 113
$ws.Range("H113").Value = 125000580
$ws.Range("I113").Value = 300
$ws.Range("K113").Value = 900
$ws.Range("M113").Value = 1270

